$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the row above (same plain style used by the new entry)
# so the new row visually matches the rest of the table.
$ws.Range("A91:J91").Copy()
$ws.Range("A97:J97").PasteSpecial(-4122)  # xlPasteFormats

# New daily COVID-19 data point for 2020-06-15 (added 2020-06-16 by the data bot)
$values = @(43997, 88165, 567, 1499, 3, 7, 1, 0, 109, 0)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(97, $i + 1).Value = $values[$i]
}

# Grow the table ("Tabela1") so it now covers the new row too
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J97"))

# Keep the current selection in sync with where Excel would leave the cursor
# after entering the new row (last cell / whole new row)
$ws.Range("A97:J97").Select() | Out-Null
